$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-07 Sunday" "2025-12-08 Monday"

Replace-Text "82×24=" "99×63="
Replace-Text "21×51=" "69×53="
Replace-Text "33×43=" "85×76="
Replace-Text "97×37=" "14×84="
Replace-Text "48×94=" "81×64="
Replace-Text "45×78=" "42×42="
Replace-Text "33×40=" "54×95="
Replace-Text "99×78=" "63×40="
Replace-Text "62×47=" "85×37="
Replace-Text "75×62=" "56×54="
Replace-Text "99×64=" "54×49="
Replace-Text "45×18=" "17×34="
Replace-Text "83×68=" "73×55="
Replace-Text "57×55=" "30×50="
Replace-Text "54×55=" "60×92="
Replace-Text "75×99=" "72×99="
Replace-Text "42×11=" "19×61="
Replace-Text "81×44=" "45×45="
Replace-Text "88×16=" "70×20="
Replace-Text "57×41=" "49×92="
Replace-Text "26×56=" "87×79="
Replace-Text "99×45=" "44×50="
Replace-Text "91×72=" "97×73="
Replace-Text "97×30=" "92×91="
Replace-Text "54×41=" "21×80="
